# Update ligand (sending-cluster) and receptor (target-cluster) TPM-derived
# statistics with newly recomputed values, then refresh the per-edge
# (Q..T) columns that are derived from them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-"Sending cluster" ligand stats: average expr (G), total expr (H)
$ligandG = @{
    "ECs"               = 144.783305
    "FAPs"              = 82.24887099999999
    "Inflammatory-Mac"  = 163.8590903333333
    "MuSCs"             = 57.0238095
    "Resolving-Mac"     = 147.8896333333333
}
$ligandH = @{
    "ECs"               = 434.349915
    "FAPs"              = 246.746613
    "Inflammatory-Mac"  = 491.577271
    "MuSCs"             = 114.047619
    "Resolving-Mac"     = 443.6689
}

# New per-"Target cluster" receptor stats: average expr (M), total expr (N)
$receptorM = @{
    "ECs"               = 0.8828746666666666
    "FAPs"              = 0.04420533333333334
    "Inflammatory-Mac"  = 4.552434333333333
    "MuSCs"             = 0.290258
    "Resolving-Mac"     = 2.847237333333334
}
$receptorN = @{
    "ECs"               = 2.648624
    "FAPs"              = 0.132616
    "Inflammatory-Mac"  = 13.657303
    "MuSCs"             = 0.580516
    "Resolving-Mac"     = 8.541712
}

# Ligand/receptor-derived specificity (I/J, O/P) are each value normalised
# against the sum across all clusters present in the sheet.
$sumG = 0
foreach ($v in $ligandG.Values) { $sumG += $v }
$sumH = 0
foreach ($v in $ligandH.Values) { $sumH += $v }
$sumM = 0
foreach ($v in $receptorM.Values) { $sumM += $v }
$sumN = 0
foreach ($v in $receptorN.Values) { $sumN += $v }

$ligandI = @{}
$ligandJ = @{}
foreach ($k in $ligandG.Keys) {
    $ligandI[$k] = $ligandG[$k] / $sumG
    $ligandJ[$k] = $ligandH[$k] / $sumH
}
$receptorO = @{}
$receptorP = @{}
foreach ($k in $receptorM.Keys) {
    $receptorO[$k] = $receptorM[$k] / $sumM
    $receptorP[$k] = $receptorN[$k] / $sumN
}

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $sending = $ws.Cells.Item($r, 1).Value2
    $target = $ws.Cells.Item($r, 4).Value2

    if (-not $ligandG.ContainsKey($sending)) { continue }
    if (-not $receptorM.ContainsKey($target)) { continue }

    $g = $ligandG[$sending]
    $h = $ligandH[$sending]
    $i = $ligandI[$sending]
    $j = $ligandJ[$sending]

    $m = $receptorM[$target]
    $n = $receptorN[$target]
    $o = $receptorO[$target]
    $p = $receptorP[$target]

    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j

    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p

    $ws.Cells.Item($r, 17).Value = $g * $m
    $ws.Cells.Item($r, 18).Value = $h * $n
    $ws.Cells.Item($r, 19).Value = $i * $o
    $ws.Cells.Item($r, 20).Value = $j * $p
}
